# Revert "more tiny changes to 1 and 2":
# Re-add the three preview bullets ("Highlighting a word", "Changing the
# theme", "Using the help tab") at the sub-bullet (lvl=1) indent level to the
# "Content Placeholder 2" body on slide 7, right after the existing
# "How to open a script" bullet and before the trailing blank paragraph.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$paraCount = $tr.Paragraphs().Count
$lastPara = $tr.Paragraphs($paraCount, 1)
[void]$lastPara.InsertAfter("`rHighlighting a word`rChanging the theme`rUsing the help tab")
